$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '60.755.57'
$ws.Range("E2").Value2 = '  +6.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '3.303.20'
$ws.Range("E3").Value2 = '  +2.03%  '

$ws.Range("E4").Value2 = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '408.75'
$ws.Range("E5").Value2 = '  +3.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '112.05'
$ws.Range("E6").Value2 = '  +4.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '3.297.65'
$ws.Range("E7").Value2 = '  +1.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.565'
$ws.Range("E8").Value2 = '  -2.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '1.00'
$ws.Range("E9").Value2 = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.616'
$ws.Range("E10").Value2 = '  -0.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.109'
$ws.Range("E11").Value2 = '  +13.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '38.48'
$ws.Range("E12").Value2 = '  -1.19%  '

$ws.Range("E13").Value2 = '  -0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '3.770.48'
$ws.Range("E14").Value2 = '  +0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '8.12'
$ws.Range("E15").Value2 = '  -0.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '18.93'
$ws.Range("E16").Value2 = '  -0.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '3.283.01'
$ws.Range("E17").Value2 = '  +0.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '60.510.31'
$ws.Range("E18").Value2 = '  +6.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '0.985'
$ws.Range("E19").Value2 = '  -4.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '10.40'
$ws.Range("E20").Value2 = '  -3.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '3.19'
$ws.Range("E22").Value2 = '  -4.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '12.31'
$ws.Range("E23").Value2 = '  -4.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '294.08'
$ws.Range("E24").Value2 = '  -0.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '73.02'
$ws.Range("E25").Value2 = '  -1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '3.05'
$ws.Range("E26").Value2 = '  -3.41%  '

$ws.Range("B27").Value2 = 'EthereumClassic'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '28.48'
$ws.Range("E27").Value2 = '  +2.40%  '

$ws.Range("B28").Value2 = 'LEO'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '4.47'
$ws.Range("E28").Value2 = '  +2.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '7.33'
$ws.Range("E29").Value2 = '  +0.34%  '

$ws.Range("E30").Value2 = '  +0.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '7.40'
$ws.Range("E31").Value2 = '  -3.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '1.00'
$ws.Range("E32").Value2 = '  +0.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '11.09'
$ws.Range("E33").Value2 = '  -2.74%  '

$ws.Range("E34").Value2 = '  -1.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '39.92'
$ws.Range("E35").Value2 = '  +4.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '2.41'
$ws.Range("E36").Value2 = '  +13.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.0474'
$ws.Range("E37").Value2 = '  -2.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '52.36'
$ws.Range("E38").Value2 = '  +1.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.998'
$ws.Range("E39").Value2 = '  -0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '3.05'
$ws.Range("E40").Value2 = '  +4.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '3.29'
$ws.Range("E41").Value2 = '  -6.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '134.81'
$ws.Range("E42").Value2 = '  -0.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.119'
$ws.Range("E43").Value2 = '  -1.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '1.87'
$ws.Range("E44").Value2 = '  -0.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '0.279'
$ws.Range("E45").Value2 = '  -0.66%  '

$ws.Range("B46").Value2 = 'NEARProtocol'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '3.72'
$ws.Range("E46").Value2 = '  -5.55%  '

$ws.Range("B47").Value2 = 'Celestia'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '16.05'
$ws.Range("E47").Value2 = '  -5.40%  '

$ws.Range("E48").Value2 = '  +3.69%  '

$ws.Range("B49").Value2 = 'EnergySwap'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '20.79'
$ws.Range("E49").Value2 = '  -5.85%  '

$ws.Range("B50").Value2 = 'Maker'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '2.106.95'
$ws.Range("E50").Value2 = '  -2.35%  '

$ws.Range("B51").Value2 = 'RocketPoolETH'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '3.657.10'
$ws.Range("E51").Value2 = '  +2.76%  '
